$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.013.45"
$ws.Range("E2").Value = "  +4.05%  "
$ws.Range("D3").Value = "2.652.21"
$ws.Range("E3").Value = "  +6.35%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "3.067.03"
$ws.Range("E15").Value = "  +6.31%  "
$ws.Range("D16").Value = "2.665.00"
$ws.Range("E16").Value = "  +7.26%  "
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "49.905.42"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "277.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.23%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.95%  "
$ws.Range("E39").Value = "  +8.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("D45").Value = "2.084.91"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.40%  "
$ws.Range("E47").Value = "  +17.15%  "
$ws.Range("E48").Value = "  +6.16%  "
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.65%  "
